$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("Use cases"): Reviewing column gets the long "Guilherme ... Leticia ... Joana" label,
# and a new "Leticia" reviewer name in the Done-adjacent Reviewing cell (E9).
$ws.Range("D9").Value = "Guilherme                 Leticia                         Joana"
$ws.Range("E9").Value = "Leticia"

# Row 10 ("Metrics"): add "Leticia" as reviewer.
$ws.Range("E10").Value = "Leticia"

# A new, empty, underlined cell was added at I9 (selection ended there).
$ws.Range("I9").Value = ""
$ws.Range("I9").Font.Underline = $true

# Page setup: paper size + portrait orientation were configured.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection ends on I9 (the last-clicked cell), and the view scrolls back to show row 1.
$null = $ws.Range("I9").Select()
